$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MILP model output values
$ws.Range("D4").Value = 0.8
$ws.Range("D6").Value = 0.6

# Update the active selection to reflect the last cell the author clicked on
$ws.Range("D13").Select()
